$d = $word.ActiveDocument

# Insert a new blank paragraph after the last existing paragraph.
$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()

# Insert a second new paragraph after that one.
$r2 = $d.Content
$r2.Collapse(0)
$r2.InsertParagraphAfter()

# Put the Test2 text into the newly added (now last) paragraph.
$r3 = $d.Content
$r3.Collapse(0)
$r3.InsertAfter("TEST2222 TEST222222 TEST22222222222222222222")
